# Insert a new header row at the top of column A and label it "ENSEMBL_ID".
# This shifts all existing ENSG... identifiers down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "ENSEMBL_ID"
